$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 72708
$ws.Range("B2").Value = 'Levi Fogaça'
$ws.Range("C2").Value = 'Jurídico'
$ws.Range("D2").Value = 'Problemas pessoais'
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45089
$ws.Range("G2").Value = 10804.19

# Row 3
$ws.Range("A3").Value = 99552
$ws.Range("B3").Value = 'André Martins'
$ws.Range("C3").Value = 'Atendimento ao Cliente'
$ws.Range("D3").Value = 'Problemas pessoais'
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45103
$ws.Range("G3").Value = 6762.29

# Row 4
$ws.Range("A4").Value = 45556
$ws.Range("B4").Value = 'Bruno Azevedo'
$ws.Range("C4").Value = 'Vendas'
$ws.Range("D4").Value = 'Outros'
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45096
$ws.Range("G4").Value = 10347.94

# Row 5
$ws.Range("A5").Value = 82670
$ws.Range("B5").Value = 'Emilly Azevedo'
$ws.Range("C5").Value = 'P&D'
$ws.Range("D5").Value = 'Viagem de negócios'
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45088
$ws.Range("G5").Value = 11259.71

# Row 6
$ws.Range("A6").Value = 87526
$ws.Range("B6").Value = 'Maria Cecília Rocha'
$ws.Range("C6").Value = 'Engenharia'
$ws.Range("D6").Value = 'Viagem de negócios'
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45093
$ws.Range("G6").Value = 3432.58

# Row 7
$ws.Range("A7").Value = 80260
$ws.Range("B7").Value = 'Guilherme Souza'
$ws.Range("C7").Value = 'P&D'
$ws.Range("D7").Value = 'Viagem de negócios'
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45081
$ws.Range("G7").Value = 9078.29

# Row 8
$ws.Range("A8").Value = 17957
$ws.Range("B8").Value = 'Kevin Barbosa'
$ws.Range("C8").Value = 'P&D'
$ws.Range("D8").Value = 'Outros'
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45105
$ws.Range("G8").Value = 5198.46

# Row 9
$ws.Range("A9").Value = 11852
$ws.Range("B9").Value = 'João Pedro Carvalho'
$ws.Range("C9").Value = 'Vendas'
$ws.Range("D9").Value = 'Outros'
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45083
$ws.Range("G9").Value = 6096.62

# Row 10
$ws.Range("A10").Value = 45738
$ws.Range("B10").Value = 'Miguel Costela'
$ws.Range("C10").Value = 'Recursos Humanos'
$ws.Range("D10").Value = 'Consulta médica'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 9994.49

# Row 11
$ws.Range("A11").Value = 78250
$ws.Range("B11").Value = 'Maria Eduarda Novaes'
$ws.Range("C11").Value = 'TI'
$ws.Range("D11").Value = 'Doença'
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45101
$ws.Range("G11").Value = 4383.67

Write-Host "applied"